# Cyclic shift of rows 20-23: the record that was on row 23 moves to row
# 20, and rows 20, 21, 22 each shift down one row (to 21, 22, 23
# respectively). Only the record-specific columns change; the shared
# columns (C, D, I, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY)
# are identical across the four rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Snapshot the "before" values for the cells that differ between the four
# records (row 23's values are captured first since row 20 will be
# overwritten with them).
$rows = @(20, 21, 22, 23)
$cols = @("A", "B", "E", "F", "G", "H", "M", "Q", "R", "Z", "AB")

$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{}
    foreach ($c in $cols) {
        $snapshot[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# New row 20 gets old row 23's data.
# New row 21 gets old row 20's data.
# New row 22 gets old row 21's data.
# New row 23 gets old row 22's data.
$sourceFor = @{ 20 = 23; 21 = 20; 22 = 21; 23 = 22 }

foreach ($destRow in $rows) {
    $srcRow = $sourceFor[$destRow]
    $src = $snapshot[$srcRow]

    $ws.Range("A$destRow").Value2 = $src["A"]
    $ws.Range("B$destRow").Value2 = $src["B"]
    $ws.Range("E$destRow").Value2 = $src["E"]
    $ws.Range("F$destRow").Value2 = $src["F"]
    $ws.Range("G$destRow").Value2 = $src["G"]
    $ws.Range("H$destRow").Value2 = $src["H"]

    # M is only populated for the "Tretåig hackspett" record; clear it on
    # rows that shouldn't have it.
    if ($src["M"] -ne $null -and $src["M"] -ne "") {
        $ws.Range("M$destRow").Value2 = $src["M"]
    } else {
        $ws.Range("M$destRow").Value2 = ""
    }

    $ws.Range("Q$destRow").Value2 = $src["Q"]
    $ws.Range("R$destRow").Value2 = $src["R"]
    $ws.Range("Z$destRow").Value2 = $src["Z"]
    $ws.Range("AB$destRow").Value2 = $src["AB"]
}
